$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): labels shifted left by one column, new "Usuario" header added in K1
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Precio"
$ws.Range("C1").Value = "Peso"
$ws.Range("D1").Value = "Tipo"
$ws.Range("E1").Value = "Contenido"
$ws.Range("F1").Value = "Categoría"
$ws.Range("G1").Value = "Dimensiones"
$ws.Range("H1").Value = "Estado pedido"
$ws.Range("I1").Value = "Direccion"
$ws.Range("J1").Value = "Domiciliario"
$ws.Range("K1").Value = "Usuario"

# Add new data row 2 with sample package values
$ws.Range("A2").Value = "taladro"
# Force text format on B2 so "10000 $" is kept as literal text, not converted to a number
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "10000 $"
$ws.Range("C2").Value = "20 kg"
$ws.Range("D2").Value = "basico"
$ws.Range("E2").Value = "taladro"
$ws.Range("F2").Value = "herramientas"
$ws.Range("G2").Value = "10x10x10"
$ws.Range("H2").Value = "pendiente"
